# scheduleData.xlsx - "Students" schedule table update.
# Several rows' time slots / durations / counts change, the data block
# shifts down starting at row 6, and two new rows (9, 10) are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (id 2): end time 08:30-11:30 -> 08:30-09:25, duration 180 -> 55
$ws.Range("C2").Value = "09:25"
$ws.Range("D2").Value = 55

# Row 3 (id 3): end time 10:00-10:28 -> 10:00-11:35, duration 28 -> 95, count 4 -> 1
$ws.Range("C3").Value = "11:35"
$ws.Range("D3").Value = 95
$ws.Range("E3").Value = 1

# Row 4 (id 4): end time 09:00-09:26 -> 09:00-09:49, duration 26 -> 49, count 2 -> 4
$ws.Range("C4").Value = "09:49"
$ws.Range("D4").Value = 49
$ws.Range("E4").Value = 4

# Row 6: now id 6, 11:00-12:15, duration 75, count 6 (new row inserted ahead of the
# old id-7 row, which moves down to row 7)
$ws.Range("A6").Value = 6
$ws.Range("B6").Value = "11:00"
$ws.Range("C6").Value = "12:15"
$ws.Range("D6").Value = 75
$ws.Range("E6").Value = 6

# Row 7: old id-7 row (10:00-10:18, 18, 2) shifted down from row 6
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "10:00"
$ws.Range("C7").Value = "10:18"
$ws.Range("D7").Value = 18
$ws.Range("E7").Value = 2

# Row 8: old id-10 row (14:00-14:50, 50, 1) shifted down from row 7
$ws.Range("A8").Value = 10
$ws.Range("B8").Value = "14:00"
$ws.Range("C8").Value = "14:50"
$ws.Range("D8").Value = 50
$ws.Range("E8").Value = 1

# Row 9 (new): old id-11 row (11:00-11:21, 21, 3) shifted down from row 8
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "11:00"
$ws.Range("C9").Value = "11:21"
$ws.Range("D9").Value = 21
$ws.Range("E9").Value = 3

# Row 10 (new): id 12, 10:00-10:26, duration 26, count 7
$ws.Range("A10").Value = 12
$ws.Range("B10").Value = "10:00"
$ws.Range("C10").Value = "10:26"
$ws.Range("D10").Value = 26
$ws.Range("E10").Value = 7
